# Final cleanup of the game stats worksheet:
#  - trim the roster down to the two "test" players (Alberto / Antonio)
#  - fix "ALberto" -> "Alberto" capitalization (folded into the rewrite below)
#  - leave a blank spacer row (row 4) above the header row
#  - park the selection back on A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old sample rows 8-17 (Ayme..The king) entirely.
$ws.Range("A8:E17").EntireRow.Delete()

# Materialize an empty spacer row at r=4 (between the title block and the
# header row) without shifting any existing rows.
$ws.Rows.Item(4).OutlineLevel = 0

# Row 6: ALberto -> Alberto, reset to the "final test" stat line.
$ws.Range("A6").Value = "Alberto"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 100

# Row 7: Iliana -> Antonio, reset to the "final test" stat line.
$ws.Range("A7").Value = "Antonio"
$ws.Range("B7").Value = 1
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# Reset the selection to A1 (was M4).
$ws.Range("A1").Select() | Out-Null
